# Updated cryptos list on Sun Feb 19 03:56:56 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row on the active sheet.
#
# Column D values are stored as plain text (e.g. "24.753.86", "1.004"),
# mirroring coinranking.com's locale-formatted numbers. Excel's normal
# Range.Value assignment auto-detects plain decimal-looking strings (like
# "1.004" or "53.91") as numbers, which would silently change the cell's
# stored type/format. To avoid that, we briefly force a Text number format
# before assigning the value, then clear the format back to the sheet's
# default (General) so no stray styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.753.86"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.699.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3922"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4032"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.501"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08872"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.217"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.036"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.704.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.009"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.749.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.244"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.43%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.93"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.163"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.735"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08734"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.071"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.170"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.22"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.968"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2737"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09153"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02731"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.461"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7659"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.81"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7157"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.572"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.212"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.59"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.308"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07975"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.78%  "
